# Gun Run Presentation - cleanup edit
#
# 1. Reorder slides: move the "Problem & future" slide (currently slide 4)
#    to position 3, ahead of the "Process" slide.
# 2. Update the title slide text from "Geometry Wars" to "Gun Run".
# 3. Consolidate the "Goal & pride" body paragraph into a single run with
#    the full (previously truncated) list of specific goals.

$p = $ppt.ActivePresentation

# --- 1. Reorder slides (swap "Process" and "Problem & future") ---
$s = $p.Slides.Item(4)
$s.MoveTo(3)

# --- 2. Title slide: "Geometry Wars" -> "Gun Run" ---
$titleSlide = $p.Slides.Item(1)
$titleShape = $titleSlide.Shapes.Item(2)
$titleShape.TextFrame.TextRange.Replace("Geometry Wars", "Gun Run", 0, 0, 0) | Out-Null

# --- 3. Goal & pride slide: merge/expand the goals paragraph ---
$goalSlide = $p.Slides.Item(2)
$goalShape = $goalSlide.Shapes.Item(2)
$goalRange = $goalShape.TextFrame.TextRange

$oldGoalText = "The goal was to create a game with a focus on challenging gameplay and interesting lore. Under this, the specific goals were"
$newGoalText = "The goal was to create a game with a focus on challenging gameplay and interesting lore. Under this, the specific goals were adding a timer, adding visual feedback for player health, a respawn mechanic for when a player falls of the world, a bullet mob, on screen text, and randomized platforms. I met all goals but the last one."

$goalRange.Replace($oldGoalText, $newGoalText, 0, 0, 0) | Out-Null
